# Update the catalog header labels on Sheet1 and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Relabel the header row: "tx" -> "t", "mean" -> "y", "SD" -> "sd" (n stays "n")
$ws.Range("D1").Value = "t"
$ws.Range("E1").Value = "n"
$ws.Range("F1").Value = "y"
$ws.Range("G1").Value = "sd"

# Move the current selection on Sheet1 to L15
$ws.Range("L15").Select()
